# Update the workbook "Översikt DALARNAS LÄN":
#  1. Column C ("Förändrad") changes from 45208 to 45212 for every data row (2-173).
#  2. For the first four data rows (2-5), the hyperlink formulas in columns
#     S, T, V, W, X and Y get an extra, more descriptive filename suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump the "Förändrad" date for every data row ------------------------
$lastRow = $ws.UsedRange.Rows.Count  # header (row 1) + 172 data rows = 173 rows total -> used range row count is 174 (includes empty row 0)
$firstDataRow = 2
$lastDataRow = 173

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # column C
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}

# --- 2. Rewrite the document-link formulas for rows 2-5 ---------------------
$baseUrl = "https://klasma.github.io/LoggingDetectiveFiles/Logging_2039"

$cases = @{
    2 = "A 30234-2023"
    3 = "A 33548-2023"
    4 = "A 33550-2023"
    5 = "A 30241-2023"
}

foreach ($row in 2..5) {
    $case = $cases[$row]

    $ws.Range("S$row").Formula = '=HYPERLINK("' + $baseUrl + '/artfynd/' + $case + ' artfynd.xlsx", "' + $case + '")'
    $ws.Range("T$row").Formula = '=HYPERLINK("' + $baseUrl + '/kartor/' + $case + ' karta.png", "' + $case + '")'
    $ws.Range("V$row").Formula = '=HYPERLINK("' + $baseUrl + '/klagomål/' + $case + ' fsc-klagomål.docx", "' + $case + '")'
    $ws.Range("W$row").Formula = '=HYPERLINK("' + $baseUrl + '/klagomålsmail/' + $case + ' fsc-klagomål mail.docx", "' + $case + '")'
    $ws.Range("X$row").Formula = '=HYPERLINK("' + $baseUrl + '/tillsyn/' + $case + ' tillsynsbegäran.docx", "' + $case + '")'
    $ws.Range("Y$row").Formula = '=HYPERLINK("' + $baseUrl + '/ti,llsynsmail/' + $case + ' tillsynsbegäran mail.docx", "' + $case + '")'
}
